$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the typo in the "Expected Outcome" cell for the DeleteUser test case
$ws.Range("E29").Value = "User is deleted"

# Update the "Outcome" cell to reflect that the test passed (matches expected outcome)
$ws.Range("F29").Value = "User is deleted"

# Add the "Justification" for the passed test case
$ws.Range("G29").Value = "Based on the given source code, this function has already been implemented"

# Move the active selection as recorded in the saved workbook
$ws.Range("H31").Select()
